$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Y24's date-time value gets switched from a date-only display to a
# full date+time display (style 3 -> style 2 in the original file).
$ws.Range("Y24").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append new row 25 with the latest bunker price data.
$ws.Range("A25").Value = 560
$ws.Range("B25").Value = 477
$ws.Range("C25").Value = 445
$ws.Range("D25").Value = 535
$ws.Range("E25").Value = 511
$ws.Range("F25").Value = 524
$ws.Range("G25").Value = 477
$ws.Range("H25").Value = 565
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 445
$ws.Range("K25").Value = 570
$ws.Range("L25").Value = 483
$ws.Range("M25").Value = 470
$ws.Range("N25").Value = 509
$ws.Range("O25").Value = 550
$ws.Range("P25").Value = 483
$ws.Range("Q25").Value = 614
$ws.Range("R25").Value = 495
$ws.Range("S25").Value = 477
$ws.Range("T25").Value = 475
$ws.Range("U25").Value = 615
$ws.Range("V25").Value = 540
$ws.Range("W25").Value = 594
$ws.Range("X25").Value = 490
$ws.Range("Y25").Value = 45755
$ws.Range("Y25").NumberFormat = "YYYY-MM-DD"
$ws.Range("Z25").Value = 846
$ws.Range("AA25").Value = 554
$ws.Range("AB25").Value = 549.5
$ws.Range("AC25").Value = 504
$ws.Range("AD25").Value = 540
$ws.Range("AE25").Value = 500
$ws.Range("AF25").Value = 502
$ws.Range("AG25").Value = 749
$ws.Range("AH25").Value = 464
$ws.Range("AI25").Value = 739
$ws.Range("AJ25").Value = 477
$ws.Range("AK25").Value = 488
$ws.Range("AL25").Value = 560
$ws.Range("AM25").Value = 540
$ws.Range("AN25").Value = 488
$ws.Range("AO25").Value = 530
$ws.Range("AP25").Value = 553
$ws.Range("AQ25").Value = 564
$ws.Range("AR25").Value = 546
$ws.Range("AS25").Value = 660
$ws.Range("AT25").Value = 630
$ws.Range("AU25").Value = 491
$ws.Range("AV25").Value = 470
